# "Specifiche di Analisi.xlsx" - add POLL/POPOLAZIONE related columns and
# split EMS/CAMS percentile columns into MIN/MAX + PASSO variants.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) -------------------------------------------------
# Existing columns keep their text but several get reshuffled, and five
# brand-new columns are inserted (D, F, I, K, L, M, N).
$ws.Range("A1").Value = "INQUINANTE"
$ws.Range("B1").Value = "ANNO"
$ws.Range("C1").Value = "PERCORSO EMS"
$ws.Range("D1").Value = "PERCORSO POLL"
$ws.Range("E1").Value = "PERCORSO GRIGLIA"
$ws.Range("F1").Value = "PERCORSO POPOLAZIONE"
$ws.Range("G1").Value = "PERCORSO KPI"
$ws.Range("H1").Value = "SCELTA GRIGLIA"
$ws.Range("I1").Value = "PERCENTILE POLL MINIMO"
$ws.Range("J1").Value = "PERCENTILE POLL MASSIMO"
$ws.Range("K1").Value = "PASSO POLL"
$ws.Range("L1").Value = "PERCENTILE EMS MINIMO"
$ws.Range("M1").Value = "PERCENTILE EMS MASSIMO"
$ws.Range("N1").Value = "PASSO EMS"
$ws.Range("O1").Value = "MAX or MEAN"
$ws.Range("P1").Value = "PERC or SOGLIA"
$ws.Range("Q1").Value = "VALORE DI SOGLIA"

# The header style (bold 16pt centered) only auto-propagates to cells that
# already existed in the sheet's old dimension (A1:L1); cells past that
# (M1:Q1) come back in default format, so restate the formatting there.
$headerRange = $ws.Range("M1:Q1")
$headerRange.Font.Bold = $true
$headerRange.Font.Size = 16
$headerRange.HorizontalAlignment = -4108

# ---- Column widths --------------------------------------------------
# The saved `width` XML attribute = ColumnWidth + 5/6, so back the target
# widths out to the ColumnWidth values that reproduce them.
$ws.Columns.Item(1).ColumnWidth = 21.276041666666668   # A
$ws.Columns.Item(2).ColumnWidth = 21.276041666666668   # B
$ws.Columns.Item(3).ColumnWidth = 27.721354166666668   # C
$ws.Columns.Item(4).ColumnWidth = 27.608072916666668   # D
$ws.Columns.Item(5).ColumnWidth = 29.053385416666668   # E
$ws.Columns.Item(6).ColumnWidth = 46.944010416666664   # F
$ws.Columns.Item(7).ColumnWidth = 29.053385416666668   # G
$ws.Columns.Item(8).ColumnWidth = 24.721354166666668   # H
$ws.Columns.Item(9).ColumnWidth = 35.498697916666664   # I
$ws.Columns.Item(10).ColumnWidth = 35.498697916666664  # J
$ws.Columns.Item(11).ColumnWidth = 19.276041666666668  # K
$ws.Columns.Item(12).ColumnWidth = 34.053385416666664  # L
$ws.Columns.Item(13).ColumnWidth = 38.721354166666664  # M
$ws.Columns.Item(14).ColumnWidth = 19.721354166666668  # N
$ws.Columns.Item(15).ColumnWidth = 21.166666666666668  # O
$ws.Columns.Item(16).ColumnWidth = 26.608072916666668  # P
$ws.Columns.Item(17).ColumnWidth = 22.276041666666668  # Q

# ---- View state -----------------------------------------------------
$ws.Activate()
$ws.Range("J1").Select()
$excel.ActiveWindow.ScrollColumn = 10
$ws.Range("O2").Select()
